# Update "想去人数" (F column) counts across sheets, as produced by the
# gh-pages data regeneration commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F6").Value = 1291
$ws.Range("F7").Value = 378
$ws.Range("F22").Value = 243
$ws.Range("F24").Value = 5413
$ws.Range("F29").Value = 359
$ws.Range("F34").Value = 39

# Sheet "演出"
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 1157
$ws.Range("F18").Value = 45

# Sheet "本地生活"
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F10").Value = 379

# Sheet "全部类型"
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F9").Value = 379
$ws.Range("F13").Value = 1291
$ws.Range("F14").Value = 378
$ws.Range("F17").Value = 1157
$ws.Range("F27").Value = 243
$ws.Range("F29").Value = 5413
$ws.Range("F35").Value = 359
$ws.Range("F49").Value = 39
